$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update harvestDate (col A) and rnaDate (col D) for rows 24-45 from 01.09.17 to 01.09.18
# Force text (not date-serial) storage, matching the original text-based cells,
# then restore the default "Normal" style so no stray number format sticks around.
$rngHarvestDate = $ws.Range("A24:A45")
$rngHarvestDate.NumberFormat = "@"
$rngHarvestDate.Value = "01.09.18"
$rngHarvestDate.Style = "Normal"

$rngRnaDate = $ws.Range("D24:D45")
$rngRnaDate.NumberFormat = "@"
$rngRnaDate.Value = "01.09.18"
$rngRnaDate.Style = "Normal"

# Update sheet view: scroll to A15, select D24:D45 with active cell D24
$null = $ws.Range("D24:D45").Select()
$ws.Application.ActiveWindow.ScrollRow = 15
